# Recursos links con titulo, importacion
#
# Inserts a new "titulo" column between "Cód. tema" and "link", renames the
# code column values to the new unit code, and fills in the titles for each
# link row. The hyperlink that used to live in column B now lives in column C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column before the old "link" column (B), pushing it to C.
#    This naturally shifts cell values/styles/column widths, matching the
#    widened dimension (A1:I6 -> A1:J6) and the column layout in the diff.
$ws.Columns("B:B").Insert()

# The inserted column doesn't inherit column A's explicit width, so copy it
# over (same width as "Cód. tema") to match the new <col min="1" max="2".../>
# run in the target layout.
$ws.Columns("B:B").ColumnWidth = $ws.Columns("A:A").ColumnWidth

# 2) The engine does not re-anchor the existing hyperlink when a column is
#    inserted, so drop the stale one (still "pointing" at B2) and recreate it
#    on the cell that now holds the URL (C2), then restore the built-in
#    hyperlink cell style so it matches the original formatting exactly.
$ws.Range("B2").Hyperlinks.Delete()
$ws.Range("C2").Hyperlinks.Add($ws.Range("C2"), "http://alimentosconstructores123.blogspot.com/") | Out-Null
$ws.Range("C2").Style = "Hipervínculo"

# 3) New header for the inserted column.
$ws.Range("B1").Value = "titulo"

# 4) New "Cód. tema" values for every data row.
$ws.Range("A2").Value = "l1-u01"
$ws.Range("A3").Value = "l1-u01"
$ws.Range("A4").Value = "l1-u01"
$ws.Range("A5").Value = "l1-u01"
$ws.Range("A6").Value = "l1-u01"

# 5) Titles for each link row.
$ws.Range("B2").Value = "titulo 1"
$ws.Range("B3").Value = "titulo 2"
$ws.Range("B4").Value = "titulo 3"
$ws.Range("B5").Value = "titulo nuevo"
$ws.Range("B6").Value = "ciencias naturales"

# 6) Match the new selection left behind in the sheet view.
$ws.Range("A4").Select()
